$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf1"
$ws.Range("C2").Value = "Cspg4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.384145666666667
$ws.Range("H2").Value = 4.152437
$ws.Range("I2").Value = 0.1014617184198512
$ws.Range("J2").Value = 0.1334061399754118
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.296849666666667
$ws.Range("N2").Value = 9.890549
$ws.Range("O2").Value = 0.06532810080989171
$ws.Range("P2").Value = 0.08920500662862836
$ws.Range("Q2").Value = 4.563320179768112
$ws.Range("R2").Value = 41.069881617913
$ws.Range("S2").Value = 0.006628301369276884
$ws.Range("T2").Value = 0.01190049560080633

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf1"
$ws.Range("C3").Value = "Cspg4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.384145666666667
$ws.Range("H3").Value = 4.152437
$ws.Range("I3").Value = 0.1014617184198512
$ws.Range("J3").Value = 0.1334061399754118
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.936479666666667
$ws.Range("N3").Value = 17.809439
$ws.Range("O3").Value = 0.1176331896601106
$ws.Range("P3").Value = 0.1606271930958688
$ws.Range("Q3").Value = 8.216952605871445
$ws.Range("R3").Value = 73.952573452843
$ws.Range("S3").Value = 0.01193526556612309
$ws.Range("T3").Value = 0.02142865380600498

# Row 4: ECs -> M1
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf1"
$ws.Range("C4").Value = "Cspg4"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.384145666666667
$ws.Range("H4").Value = 4.152437
$ws.Range("I4").Value = 0.1014617184198512
$ws.Range("J4").Value = 0.1334061399754118
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.3593903333333333
$ws.Range("N4").Value = 1.078171
$ws.Range("O4").Value = 0.00712143115395331
$ws.Range("P4").Value = 0.00972425809748224
$ws.Range("Q4").Value = 0.4974485725252222
$ws.Range("R4").Value = 4.477037152727
$ws.Range("S4").Value = 0.0007225526424887666
$ws.Range("T4").Value = 0.001297275736909747

# Row 5: ECs -> M2
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Fgf1"
$ws.Range("C5").Value = "Cspg4"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.384145666666667
$ws.Range("H5").Value = 4.152437
$ws.Range("I5").Value = 0.1014617184198512
$ws.Range("J5").Value = 0.1334061399754118
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3495983333333333
$ws.Range("N5").Value = 1.048795
$ws.Range("O5").Value = 0.006927399630587783
$ws.Range("P5").Value = 0.009459309582013322
$ws.Range("Q5").Value = 0.4838950181572222
$ws.Range("R5").Value = 4.355055163415
$ws.Range("S5").Value = 0.0007028658707004787
$ws.Range("T5").Value = 0.001261929978168823

# Row 6: ECs -> sCs
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Fgf1"
$ws.Range("C6").Value = "Cspg4"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.384145666666667
$ws.Range("H6").Value = 4.152437
$ws.Range("I6").Value = 0.1014617184198512
$ws.Range("J6").Value = 0.1334061399754118
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 40.5237085
$ws.Range("N6").Value = 81.047417
$ws.Range("O6").Value = 0.8029898787454566
$ws.Range("P6").Value = 0.7309842325960072
$ws.Range("Q6").Value = 56.09071551753816
$ws.Range("R6").Value = 336.544293105229
$ws.Range("S6").Value = 0.08147273297126195
$ws.Range("T6").Value = 0.0975177848535219

# Row 7: FAPs -> ECs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fgf1"
$ws.Range("C7").Value = "Cspg4"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.458038666666667
$ws.Range("H7").Value = 7.374116000000001
$ws.Range("I7").Value = 0.1801810554109116
$ws.Range("J7").Value = 0.2369096391566985
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.296849666666667
$ws.Range("N7").Value = 9.890549
$ws.Range("O7").Value = 0.06532810080989171
$ws.Range("P7").Value = 0.08920500662862836
$ws.Range("Q7").Value = 8.103783958853779
$ws.Range("R7").Value = 72.93405562968401
$ws.Range("S7").Value = 0.01177088615191671
$ws.Range("T7").Value = 0.02113352593135924

# Row 8: FAPs -> FAPs
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Fgf1"
$ws.Range("C8").Value = "Cspg4"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.458038666666667
$ws.Range("H8").Value = 7.374116000000001
$ws.Range("I8").Value = 0.1801810554109116
$ws.Range("J8").Value = 0.2369096391566985
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 5.936479666666667
$ws.Range("N8").Value = 17.809439
$ws.Range("O8").Value = 0.1176331896601106
$ws.Range("P8").Value = 0.1606271930958688
$ws.Range("Q8").Value = 14.59209656454711
$ws.Range("R8").Value = 131.328869080924
$ws.Range("S8").Value = 0.02119527226431065
$ws.Range("T8").Value = 0.03805413035509563

# Row 9: FAPs -> M1
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Fgf1"
$ws.Range("C9").Value = "Cspg4"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.458038666666667
$ws.Range("H9").Value = 7.374116000000001
$ws.Range("I9").Value = 0.1801810554109116
$ws.Range("J9").Value = 0.2369096391566985
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.3593903333333333
$ws.Range("N9").Value = 1.078171
$ws.Range("O9").Value = 0.00712143115395331
$ws.Range("P9").Value = 0.00972425809748224
$ws.Range("Q9").Value = 0.8833953357595556
$ws.Range("R9").Value = 7.950558021836001
$ws.Range("S9").Value = 0.001283146981355453
$ws.Range("T9").Value = 0.002303770476941121

# Row 10: FAPs -> M2
$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Fgf1"
$ws.Range("C10").Value = "Cspg4"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.458038666666667
$ws.Range("H10").Value = 7.374116000000001
$ws.Range("I10").Value = 0.1801810554109116
$ws.Range("J10").Value = 0.2369096391566985
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.3495983333333333
$ws.Range("N10").Value = 1.048795
$ws.Range("O10").Value = 0.006927399630587783
$ws.Range("P10").Value = 0.009459309582013322
$ws.Range("Q10").Value = 0.8593262211355556
$ws.Range("R10").Value = 7.73393599022
$ws.Range("S10").Value = 0.001248186176692466
$ws.Range("T10").Value = 0.002241001619746277

# Row 11: FAPs -> sCs
$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Fgf1"
$ws.Range("C11").Value = "Cspg4"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 2.458038666666667
$ws.Range("H11").Value = 7.374116000000001
$ws.Range("I11").Value = 0.1801810554109116
$ws.Range("J11").Value = 0.2369096391566985
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 40.5237085
$ws.Range("N11").Value = 81.047417
$ws.Range("O11").Value = 0.8029898787454566
$ws.Range("P11").Value = 0.7309842325960072
$ws.Range("Q11").Value = 99.60884240972867
$ws.Range("R11").Value = 597.653054458372
$ws.Range("S11").Value = 0.1446835638366363
$ws.Range("T11").Value = 0.1731772107735562

# Row 12: sCs -> ECs
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Fgf1"
$ws.Range("C12").Value = "Cspg4"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 9.799864
$ws.Range("H12").Value = 19.599728
$ws.Range("I12").Value = 0.7183572261692373
$ws.Range("J12").Value = 0.6296842208678898
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 3.296849666666667
$ws.Range("N12").Value = 9.890549
$ws.Range("O12").Value = 0.06532810080989171
$ws.Range("P12").Value = 0.08920500662862836
$ws.Range("Q12").Value = 32.30867836177867
$ws.Range("R12").Value = 193.852070170672
$ws.Range("S12").Value = 0.04692891328869811
$ws.Range("T12").Value = 0.0561709850964628

# Row 13: sCs -> FAPs
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Fgf1"
$ws.Range("C13").Value = "Cspg4"
$ws.Range("D13").Value = "FAPs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 9.799864
$ws.Range("H13").Value = 19.599728
$ws.Range("I13").Value = 0.7183572261692373
$ws.Range("J13").Value = 0.6296842208678898
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 5.936479666666667
$ws.Range("N13").Value = 17.809439
$ws.Range("O13").Value = 0.1176331896601106
$ws.Range("P13").Value = 0.1606271930958688
$ws.Range("Q13").Value = 58.17669337209867
$ws.Range("R13").Value = 349.060160232592
$ws.Range("S13").Value = 0.08450265182967684
$ws.Range("T13").Value = 0.1011444089347683

# Row 14: sCs -> M1
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Fgf1"
$ws.Range("C14").Value = "Cspg4"
$ws.Range("D14").Value = "M1"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 9.799864
$ws.Range("H14").Value = 19.599728
$ws.Range("I14").Value = 0.7183572261692373
$ws.Range("J14").Value = 0.6296842208678898
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.3593903333333333
$ws.Range("N14").Value = 1.078171
$ws.Range("O14").Value = 0.00712143115395331
$ws.Range("P14").Value = 0.00972425809748224
$ws.Range("Q14").Value = 3.521976389581333
$ws.Range("R14").Value = 21.131858337488
$ws.Range("S14").Value = 0.00511573153010909
$ws.Range("T14").Value = 0.006123211883631373

# Row 15: sCs -> M2
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Fgf1"
$ws.Range("C15").Value = "Cspg4"
$ws.Range("D15").Value = "M2"
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 9.799864
$ws.Range("H15").Value = 19.599728
$ws.Range("I15").Value = 0.7183572261692373
$ws.Range("J15").Value = 0.6296842208678898
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0.3495983333333333
$ws.Range("N15").Value = 1.048795
$ws.Range("O15").Value = 0.006927399630587783
$ws.Range("P15").Value = 0.009459309582013322
$ws.Range("Q15").Value = 3.426016121293333
$ws.Range("R15").Value = 20.55609672776
$ws.Range("S15").Value = 0.004976347583194839
$ws.Range("T15").Value = 0.005956377984098223

# Row 16: sCs -> sCs
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Fgf1"
$ws.Range("C16").Value = "Cspg4"
$ws.Range("D16").Value = "sCs"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 9.799864
$ws.Range("H16").Value = 19.599728
$ws.Range("I16").Value = 0.7183572261692373
$ws.Range("J16").Value = 0.6296842208678898
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 40.5237085
$ws.Range("N16").Value = 81.047417
$ws.Range("O16").Value = 0.8029898787454566
$ws.Range("P16").Value = 0.7309842325960072
$ws.Range("Q16").Value = 397.126832075644
$ws.Range("R16").Value = 1588.507328302576
$ws.Range("S16").Value = 0.5768335819375583
$ws.Range("T16").Value = 0.4602892369689291
